$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10, shifting existing rows 10-27 down to 11-28
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price record
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = 44481
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 300000000
$ws.Cells.Item(10, 7).Value = "Espárragos"
$ws.Cells.Item(10, 8).Value = "Verde"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 4000
$ws.Cells.Item(10, 11).Value = 900
$ws.Cells.Item(10, 12).Value = 900
$ws.Cells.Item(10, 13).Value = 900
$ws.Cells.Item(10, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Linares"
$ws.Cells.Item(10, 16).Value = 90
$ws.Cells.Item(10, 17).Value = 10
$ws.Cells.Item(10, 18).Value = "Hortaliza"
